# Agronomy - Irrigation template: fix the Input/Output column headers on the
# "Events-Irrigation" sheet (and its annotationTable) so they reference
# ISA "Sample Name" instead of the erroneous "Source Name" - resolving the
# i/o naming errors called out in the commit message.

$wb = $excel.ActiveWorkbook

# The data/annotation table lives on the second sheet.
$ws = $wb.Worksheets.Item("Events-Irrigation")

# Column A header: "Input [Source Name]" -> "Input [Sample Name]"
$ws.Range("A1").Value = "Input [Sample Name]"

# Column V header: "Output [Source Name]" -> "Output [Sample Name]"
$ws.Range("V1").Value = "Output [Sample Name]"
